# Week 13 logging update - Target Depth Data
$wb = $excel.ActiveWorkbook

# OFF sheet - Home (H) row updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 405
$wsOff.Range("C2").Value = 302
$wsOff.Range("D2").Value = 118
$wsOff.Range("E2").Value = 62

# DEF sheet - Home (H) row updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 443
$wsDef.Range("C2").Value = 307
$wsDef.Range("D2").Value = 111
$wsDef.Range("E2").Value = 43
$wsDef.Range("G2").Value = 4
